# Branch - Brand - SKU wise Stock Aging Status: Summary -- data refresh.
# The underlying SKU ordering/aging data was refreshed; this reshuffles a
# handful of same-brand SKU rows (Item Name / UOM swap in place) and
# updates the refreshed numeric figures (TP Sales Value / stock / sales
# columns) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dinafex (rows 3-5): rotate item order, TP Sales Value (BB) follows ---
$ws.Range("C3").Value = "Dinafex 120mg Tablet"
$ws.Range("D3").Value = "30's"
$ws.Range("BB3").Value = 179.91

$ws.Range("C4").Value = "Dinafex 60mg Tablet"
$ws.Range("D4").Value = "30's"
$ws.Range("BB4").Value = 78.70999999999999

$ws.Range("C5").Value = "Dinafex 180mg Tablet"
$ws.Range("D5").Value = "30's"
$ws.Range("BB5").Value = 224.89

# --- Etorix (rows 7-9): rotate item order, TP Sales Value (BB) follows ---
$ws.Range("C7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D7").Value = "40's"

$ws.Range("C8").Value = "Etorix 90mg Tablet"
$ws.Range("D8").Value = "30's"
$ws.Range("BB8").Value = 269.87

$ws.Range("C9").Value = "Etorix 120mg Tablet"
$ws.Range("D9").Value = "20's"
$ws.Range("BB9").Value = 209.9

# --- Geminox 320mg Tablet - 8's (row 13): Monthly Sales Target updated ---
$ws.Range("H13").Value = 217

# --- Ketonic (rows 15-16): swap item order, TP Sales Value (BB) follows ---
$ws.Range("C15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("D15").Value = "4's"
$ws.Range("BB15").Value = 165.41

$ws.Range("C16").Value = "Ketonic 10mg Tablet"
$ws.Range("D16").Value = "20's"
$ws.Range("BB16").Value = 150.38

# --- Kynol (rows 18-19): swap item order, TP Sales Value (BB) follows ---
$ws.Range("C18").Value = "Kynol TR 100mg Capsule"
$ws.Range("D18").Value = "50 's"
$ws.Range("BB18").Value = 262.37

$ws.Range("C19").Value = "Kynol D 25mg Tablet"
$ws.Range("D19").Value = "60 's"
$ws.Range("BB19").Value = 180.45

# --- Flucloxin 500mg Capsule - 36's (row 23): refreshed sales/stock figures ---
$ws.Range("H23").Value = 1020
$ws.Range("L23").Value = 14
$ws.Range("M23").Value = 0

# --- Sk-Mox 500mg Capsule (row 24): refreshed sales/stock figures ---
$ws.Range("E24").Value = 0
$ws.Range("I24").Value = 12
$ws.Range("L24").Value = 21
$ws.Range("N24").Value = 163
$ws.Range("O24").Value = 175
$ws.Range("T24").Value = 175
$ws.Range("AU24").Value = 154
$ws.Range("AZ24").Value = 131
$ws.Range("BA24").Value = 28706
$ws.Range("BC24").Value = 0
$ws.Range("BD24").Value = 0

# --- Zithrox (rows 25-28): rotate item order, TP Sales Value (BB) follows ---
$ws.Range("C25").Value = "Zithrox 500mg Tablet"
$ws.Range("D25").Value = "6 's"
$ws.Range("BB25").Value = 136.83

$ws.Range("C26").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D26").Value = "30ml"
$ws.Range("BB26").Value = 97.45

$ws.Range("C27").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("D27").Value = "6's"
$ws.Range("BB27").Value = 89.95999999999999

$ws.Range("C28").Value = "Zithrox 15ml Suspension"
$ws.Range("D28").Value = "15 ml"
$ws.Range("BB28").Value = 71.95999999999999
